$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows 2-3 down to 3-4
$ws.Rows("2:2").Insert()
$ws.Rows("2:2").ClearFormats()

# Re-point the dimension / header handled automatically by cell writes below

# --- Header row 1: move Odd_CS_3-3_HT from BC to AW, shifting AW..BC right by one ---
$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Time"
$ws.Range("D1").Value = "League"
$ws.Range("E1").Value = "Home"
$ws.Range("F1").Value = "Away"
$ws.Range("G1").Value = "Odd_H_FT"
$ws.Range("H1").Value = "Odd_D_FT"
$ws.Range("I1").Value = "Odd_A_FT"
$ws.Range("J1").Value = "Odd_H_HT"
$ws.Range("K1").Value = "Odd_D_HT"
$ws.Range("L1").Value = "Odd_A_HT"
$ws.Range("M1").Value = "Odd_Over05_FT"
$ws.Range("N1").Value = "Odd_Under05_FT"
$ws.Range("O1").Value = "Odd_Over15_FT"
$ws.Range("P1").Value = "Odd_Under15_FT"
$ws.Range("Q1").Value = "Odd_Over25_FT"
$ws.Range("R1").Value = "Odd_Under25_FT"
$ws.Range("S1").Value = "Odd_Over05_HT"
$ws.Range("T1").Value = "Odd_Under05_HT"
$ws.Range("U1").Value = "Odd_BTTS_Yes"
$ws.Range("V1").Value = "Odd_BTTS_No"
$ws.Range("W1").Value = "Odd_CS_1-0"
$ws.Range("X1").Value = "Odd_CS_2-0"
$ws.Range("Y1").Value = "Odd_CS_2-1"
$ws.Range("Z1").Value = "Odd_CS_3-0"
$ws.Range("AA1").Value = "Odd_CS_3-1"
$ws.Range("AB1").Value = "Odd_CS_3-2"
$ws.Range("AC1").Value = "Odd_CS_0-0"
$ws.Range("AD1").Value = "Odd_CS_1-1"
$ws.Range("AE1").Value = "Odd_CS_2-2"
$ws.Range("AF1").Value = "Odd_CS_3-3"
$ws.Range("AG1").Value = "Odd_CS_4-4"
$ws.Range("AH1").Value = "Odd_CS_0-1"
$ws.Range("AI1").Value = "Odd_CS_0-2"
$ws.Range("AJ1").Value = "Odd_CS_1-2"
$ws.Range("AK1").Value = "Odd_CS_0-3"
$ws.Range("AL1").Value = "Odd_CS_1-3"
$ws.Range("AM1").Value = "Odd_CS_2-3"
$ws.Range("AN1").Value = "Odd_CS_1-0_HT"
$ws.Range("AO1").Value = "Odd_CS_2-0_HT"
$ws.Range("AP1").Value = "Odd_CS_2-1_HT"
$ws.Range("AQ1").Value = "Odd_CS_3-0_HT"
$ws.Range("AR1").Value = "Odd_CS_3-1_HT"
$ws.Range("AS1").Value = "Odd_CS_3-2_HT"
$ws.Range("AT1").Value = "Odd_CS_0-0_HT"
$ws.Range("AU1").Value = "Odd_CS_1-1_HT"
$ws.Range("AV1").Value = "Odd_CS_2-2_HT"
$ws.Range("AW1").Value = "Odd_CS_3-3_HT"
$ws.Range("AX1").Value = "Odd_CS_0-1_HT"
$ws.Range("AY1").Value = "Odd_CS_0-2_HT"
$ws.Range("AZ1").Value = "Odd_CS_1-2_HT"
$ws.Range("BA1").Value = "Odd_CS_0-3_HT"
$ws.Range("BB1").Value = "Odd_CS_1-3_HT"
$ws.Range("BC1").Value = "Odd_CS_2-3_HT"
$ws.Range("BD1").Value = "Odd_CS_4-4_HT"

# --- Row 2 data ---
$ws.Range("A2").Value = "MgX5nwjB"
$ws.Range("B2").Value = "27/11/2024"
$ws.Range("C2").Value = "16:00"
$ws.Range("D2").Value = "BOLIVIA - DIVISION PROFESIONAL"
$ws.Range("E2").Value = "GV San Jose"
$ws.Range("F2").Value = "Bolivar"
$ws.Range("G2").Value = 3.6
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 1.9
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 2.4
$ws.Range("L2").Value = 2.4
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 17
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("Q2").Value = 1.5
$ws.Range("R2").Value = 2.5
$ws.Range("S2").Value = 1.29
$ws.Range("T2").Value = 3.5
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 17
$ws.Range("X2").Value = 21
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 41
$ws.Range("AA2").Value = 26
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 8
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 101
$ws.Range("AH2").Value = 11
$ws.Range("AI2").Value = 12
$ws.Range("AJ2").Value = 9
$ws.Range("AK2").Value = 17
$ws.Range("AL2").Value = 13
$ws.Range("AM2").Value = 19
$ws.Range("AN2").Value = 6
$ws.Range("AO2").Value = 19
$ws.Range("AP2").Value = 21
$ws.Range("AQ2").Value = 51
$ws.Range("AR2").Value = 51
$ws.Range("AS2").Value = 101
$ws.Range("AT2").Value = 3.5
$ws.Range("AU2").Value = 7
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 301
$ws.Range("AX2").Value = 4.33
$ws.Range("AY2").Value = 10
$ws.Range("AZ2").Value = 17
$ws.Range("BA2").Value = 29
$ws.Range("BB2").Value = 41
$ws.Range("BC2").Value = 81
$ws.Range("BD2").ClearContents()

# --- Row 3 data ---
$ws.Range("A3").Value = "C66nTKo1"
$ws.Range("B3").Value = "27/11/2024"
$ws.Range("C3").Value = "15:00"
$ws.Range("D3").Value = "SPAIN - LALIGA2"
$ws.Range("E3").Value = "Castellon"
$ws.Range("F3").Value = "Racing Club Ferrol"
$ws.Range("G3").Value = 1.7
$ws.Range("H3").Value = 3.7
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 2.3
$ws.Range("K3").Value = 2.3
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 12
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 1.36
$ws.Range("T3").Value = 3
$ws.Range("U3").Value = 1.73
$ws.Range("V3").Value = 2
$ws.Range("W3").Value = 7.5
$ws.Range("X3").Value = 8.5
$ws.Range("Y3").Value = 8.5
$ws.Range("Z3").Value = 13
$ws.Range("AA3").Value = 13
$ws.Range("AB3").Value = 23
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 7
$ws.Range("AE3").Value = 15
$ws.Range("AF3").Value = 41
$ws.Range("AG3").Value = 201
$ws.Range("AH3").Value = 15
$ws.Range("AI3").Value = 26
$ws.Range("AJ3").Value = 15
$ws.Range("AK3").Value = 51
$ws.Range("AL3").Value = 41
$ws.Range("AM3").Value = 41
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 8.5
$ws.Range("AP3").Value = 19
$ws.Range("AQ3").Value = 26
$ws.Range("AR3").Value = 41
$ws.Range("AS3").Value = 126
$ws.Range("AT3").Value = 3.25
$ws.Range("AU3").Value = 8
$ws.Range("AV3").Value = 51
$ws.Range("AW3").Value = 81
$ws.Range("AX3").Value = 6.5
$ws.Range("AY3").Value = 26
$ws.Range("AZ3").Value = 29
$ws.Range("BA3").Value = 81
$ws.Range("BB3").Value = 101
$ws.Range("BC3").Value = 201
$ws.Range("BD3").Value = 81

# --- Row 4 data ---
$ws.Range("A4").Value = "KCTDqtWs"
$ws.Range("B4").Value = "27/11/2024"
$ws.Range("C4").Value = "15:00"
$ws.Range("D4").Value = "SPAIN - LALIGA2"
$ws.Range("E4").Value = "Levante"
$ws.Range("F4").Value = "Malaga"
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.38
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.88
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.75
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.83
$ws.Range("W4").Value = 7
$ws.Range("X4").Value = 8
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 13
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 10
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 251
$ws.Range("AH4").Value = 13
$ws.Range("AI4").Value = 26
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 51
$ws.Range("AL4").Value = 41
$ws.Range("AM4").Value = 41
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 9
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 29
$ws.Range("AR4").Value = 51
$ws.Range("AS4").Value = 151
$ws.Range("AT4").Value = 2.75
$ws.Range("AU4").Value = 8.5
$ws.Range("AV4").Value = 51
$ws.Range("AW4").Value = 81
$ws.Range("AX4").Value = 6.5
$ws.Range("AY4").Value = 26
$ws.Range("AZ4").Value = 34
$ws.Range("BA4").Value = 81
$ws.Range("BB4").Value = 101
$ws.Range("BC4").Value = 251
$ws.Range("BD4").Value = 81
